$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final data grid (A1:G4). A new expense record ("Various paper supplies")
# was inserted at the top of the table, and the existing three records were
# reshuffled/re-columned so that the date now lives in column G instead of
# column C (amount/currency columns shift left to fill C:F).
$data = @(
    @("339040", "Various paper supplies",      "82424",  "16484.8", "98908.8", "EUR", "2017-02-17"),
    @("504240", "IT Support",                  "260053", "52010.6", "312064",  "CAD", "2017-04-26"),
    @("254423", "Waste management services",   "97485",  "19497",   "116982",  "CAD", "2017-08-15"),
    @("614399", "Concierge Services",           "279738", "55947.6", "335686",  "USD", "2017-08-17")
)

# Every cell in this sheet is stored as plain text (shared string), even the
# numeric-looking IDs/amounts and the ISO dates -- there is no header row and
# no special number formatting anywhere (every cell uses the workbook's
# default style). Writing the raw strings straight into `.Value` would let
# Excel auto-coerce "339040"/"2017-02-17" into a real number/date (and stamp
# a date NumberFormat on that cell), so we force text storage with the
# NumberFormat "@" trick and then restore the original default cell style by
# pasting the formats from an untouched, never-written cell back on top -
# this keeps every written cell on the workbook's existing default style
# instead of leaving it on the temporary "@" style.
for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $cell = $ws.Cells.Item($r + 1, $c + 1)
        $cell.NumberFormat = "@"
        $cell.Value = $row[$c]
    }
}

$blankDefaultStyleCell = $ws.Cells.Item(50, 50)
$blankDefaultStyleCell.Copy()
$ws.Range("A1:G4").PasteSpecial(-4122)
